# Generate Report for Handoff
# - Flips the localization "Status" from "In Translation" to "Ready for
#   handoff" everywhere it appears (Overview zh-cn/de-de status columns,
#   and the per-language "Status" column on the zh-cn / de-de sheets).
# - Refreshes the "Latest Handoff"/"Latest HO Xliff Generate Date"
#   timestamps to reflect the new handoff generation time.
# - Widens the now-longer Status-ish columns so the new text fits
#   (mirrors Excel's own column autosize after the text grew).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ----------------------
$overview.Range("E2").Value = "Ready for handoff"   # zh-cn status column
$overview.Range("F2").Value = "Ready for handoff"   # de-de status column
$zhcn.Range("C2").Value     = "Ready for handoff"   # Status
$dede.Range("C2").Value     = "Ready for handoff"   # Status

# --- Timestamps -------------------------------------------------------------
$overview.Range("G2").Value = "2016-09-03 19:12:23"  # Latest HO Xliff Generate Date
$dede.Range("H2").Value     = "2016-09-03 19:12:23"  # Latest Handoff Datetime (de-de)
$zhcn.Range("H2").Value     = "2016-09-03 19:12:19"  # Latest Handoff Datetime (zh-cn)

# --- Column widths: autofit the widened Status columns ---------------------
$overview.Columns.Item(5).ColumnWidth = 16.3   # E:E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 16.3   # F:F (de-de status)
$zhcn.Columns.Item(3).ColumnWidth     = 16.3   # C:C (Status)
$dede.Columns.Item(3).ColumnWidth     = 16.3   # C:C (Status)
